{"js": "// Add \"Version 4.1.7\", \"Version 4.1.8\", \"Version 4.1.9\" entries to the\n// bulleted version-notes list, continuing directly after the last existing\n// bullet (\"Moved stat page URLs to top of script, will make for easier\n// updating\"), matching the numbering (numId 15) / indent-level pattern\n// already used for the other \"Version x.y.z\" headings and their\n// sub-bullets.\n\nconst RSQUOTE = \"\\u2019\"; // U+2019 RIGHT SINGLE QUOTATION MARK\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Anchor on the very last paragraph in the document body (the\n// \"Moved stat page URLs...\" bullet that precedes the new content).\nlet anchor = paragraphs.items[paragraphs.items.length - 1];\n\nfunction addBullet(text, level) {\n  anchor = anchor.insertParagraph(text, \"After\");\n  anchor.listItem.level = level;\n  return anchor;\n}\n\n// Version 4.1.7\naddBullet(\"Version 4.1.7\", 2);\naddBullet(\n  \"added a Pythagorean win stat, using formula from Football Outsiders\",\n  3\n);\naddBullet(\n  \"only included for final 2021 week (Week 15), will not be included in whole season, will be included for whole of future seasons\",\n  4\n);\naddBullet(\n  \"fixed error in points per play stat where points per game was being used for calculation instead of total points\",\n  3\n);\naddBullet(\n  \"this probably didn\" +\n    RSQUOTE +\n    \"t affect the rankings too much, but may have negatively impacted teams which played an extra game (conference championship, maybe some bowl games)\",\n  4\n);\n\n// Version 4.1.8\naddBullet(\"Version 4.1.8\", 2);\naddBullet(\n  \"added adjusted net yards per passing attempt stat, using ESPN\" +\n    RSQUOTE +\n    \"s Bill Connelly\" +\n    RSQUOTE +\n    \"s formula\",\n  3\n);\naddBullet(\n  \"net yards per attempt, plus 20 yards per TD and -45 per INT\",\n  4\n);\n\n// Version 4.1.9\naddBullet(\"Version 4.1.9\", 2);\naddBullet(\n  \"ESPN updated their team names back to how it had been for all but the last fucking week of the season, so I had to too, meaning all team names include school name plus nickname\",\n  3\n);\n\nawait context.sync();\n", "ps1": "# Add \"Version 4.1.7\", \"Version 4.1.8\", \"Version 4.1.9\" entries to the\n# bulleted version-notes list, continuing directly after the last existing\n# bullet (\"Moved stat page URLs to top of script, will make for easier\n# updating\"), matching the numbering / indent-level pattern already used\n# for the other \"Version x.y.z\" headings and their sub-bullets.\n\n$RSQUOTE = [char]0x2019   # U+2019 RIGHT SINGLE QUOTATION MARK\n\n$d = $word.ActiveDocument\n\nfunction Add-Bullet($text, $level) {\n    $last = $d.Paragraphs.Last\n    $last.Range.InsertParagraphAfter()\n    $newPara = $d.Paragraphs.Last\n    $newPara.Range.Text = $text\n    $newPara.Range.ListFormat.ListLevelNumber = $level\n    return $newPara\n}\n\n# Version 4.1.7\nAdd-Bullet \"Version 4.1.7\" 3\nAdd-Bullet \"added a Pythagorean win stat, using formula from Football Outsiders\" 4\nAdd-Bullet \"only included for final 2021 week (Week 15), will not be included in whole season, will be included for whole of future seasons\" 5\nAdd-Bullet \"fixed error in points per play stat where points per game was being used for calculation instead of total points\" 4\nAdd-Bullet (\"this probably didn\" + $RSQUOTE + \"t affect the rankings too much, but may have negatively impacted teams which played an extra game (conference championship, maybe some bowl games)\") 5\n\n# Version 4.1.8\nAdd-Bullet \"Version 4.1.8\" 3\nAdd-Bullet (\"added adjusted net yards per passing attempt stat, using ESPN\" + $RSQUOTE + \"s Bill Connelly\" + $RSQUOTE + \"s formula\") 4\nAdd-Bullet \"net yards per attempt, plus 20 yards per TD and -45 per INT\" 5\n\n# Version 4.1.9\nAdd-Bullet \"Version 4.1.9\" 3\nAdd-Bullet \"ESPN updated their team names back to how it had been for all but the last fucking week of the season, so I had to too, meaning all team names include school name plus nickname\" 4\n"}
